# FeScienceTimingLogBook2024-2025.xlsx edit
# "moved all to first Brilluoin zone. Background removal for 400meV and various file transformations"
#
# Fills in the "D" (notes) column for rows 36-60 with task annotations,
# changes a couple of existing "WE" entries in column B to "EH", adds
# missing B entries for rows 64-65, and updates the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Set the new text values in column D (and the two new B cells)
# ---------------------------------------------------------------
$ws.Range("D36").Value = "magneticFF/sqw_op+bin"
$ws.Range("D37").Value = "backgroun&Symmetry calc"
$ws.Range("D40").Value = "PH"
$ws.Range("D41").Value = "pixel cahce"
$ws.Range("D42").Value = "pixel cahce/plotting review"
$ws.Range("D43").Value = "-combine in sqw_op"
$ws.Range("D44").Value = "-combine in sqw_op"
$ws.Range("D45").Value = "WE"
$ws.Range("D46").Value = "WE"
$ws.Range("D47").Value = "-combine in sqw_op"
$ws.Range("D48").Value = "-combine in sqw_op"
$ws.Range("D49").Value = "-combine in sqw_op"
$ws.Range("D50").Value = "background/symmetry calc"
$ws.Range("D51").Value = "combine/test interdependencies "
$ws.Range("D52").Value = "WE"
$ws.Range("D53").Value = "WE"
$ws.Range("D59").Value = "WE"
$ws.Range("D60").Value = "WE"

$ws.Range("B64").Value = "sqw_op_bin_pixels"
$ws.Range("B65").Value = "magneticFF"

# ---------------------------------------------------------------
# 2. Re-point two existing "WE" entries in column B to "EH"
# ---------------------------------------------------------------
$ws.Range("B53").Value = "EH"
$ws.Range("B56").Value = "EH"

# ---------------------------------------------------------------
# 3. Apply the correct cell formatting by copying it from existing
#    donor cells that already carry the required style, so the
#    workbook's style table / indices are reused rather than
#    duplicated. PasteSpecial only honours the first area of a
#    multi-area range, so paste into each target cell individually.
# ---------------------------------------------------------------
$xlPasteFormats = -4122

# Style used by D38 / D39 / B40 / B41 ... (red "Bad" highlight)
$redCells = @("D40","D45","D46","D52","D53","D59","D60")
$ws.Range("D38").Copy() | Out-Null
foreach ($addr in $redCells) {
    $ws.Range($addr).PasteSpecial($xlPasteFormats) | Out-Null
}

# Style used by B44 / B58 (date-number-format + blue fill)
$blueCells = @("D37","D50")
$ws.Range("B44").Copy() | Out-Null
foreach ($addr in $blueCells) {
    $ws.Range($addr).PasteSpecial($xlPasteFormats) | Out-Null
}

# Style used by B39 / B42 / B43 (plain, quote-prefixed)
$plainQuoteCells = @("D36","D43","D44","D47","D48","D49")
$ws.Range("B39").Copy() | Out-Null
foreach ($addr in $plainQuoteCells) {
    $ws.Range($addr).PasteSpecial($xlPasteFormats) | Out-Null
}

$excel.CutCopyMode = $false

# ---------------------------------------------------------------
# 4. Update the view: scroll position and active selection
# ---------------------------------------------------------------
$ws.Activate()
$ws.Range("D52").Select()
$excel.ActiveWindow.ScrollRow = 35
$excel.ActiveWindow.ScrollColumn = 1

Write-Output "Edit complete"
